$wb = $excel.ActiveWorkbook
$win = $excel.ActiveWindow
Write-Host ($win | Get-Member | Out-String)
